# Update the END-USER upload sheet: replace the 5 existing rows of
# name/email/employee-code data with 3 new rows, dropping the last two
# rows (5 and 6) entirely and refreshing the mailto hyperlinks on column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the old hyperlinks (and their relationships) up front ---
$ws.Hyperlinks.Delete()

# --- 2. Delete rows 5 and 6 (bottom-up so indices stay valid) ---
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()

# --- 3. Write the new NAME / EMAIL / EMPLOYEE_CODE values ---
$ws.Range("A2").Value = "Jayashree Kulai"
$ws.Range("B2").Value = "jayashree.cs16@sahyadri.edu.in"
$ws.Range("C2").Value = "MNG001"

$ws.Range("A3").Value = "Roy Pashan"
$ws.Range("B3").Value = "mail2winstonroy@yahoo.com"
$ws.Range("C3").Value = "MNG002"

$ws.Range("A4").Value = "Winston"
$ws.Range("B4").Value = "pashanwinsty1998@gmail.com"
$ws.Range("C4").Value = "ADM001"

# --- 4. Re-create the mailto hyperlinks on the email column ---
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:jayashree.cs16@sahyadri.edu.in", "", "", "jayashree.cs16@sahyadri.edu.in")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:mail2winstonroy@yahoo.com", "", "", "mail2winstonroy@yahoo.com")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:pashanwinsty1998@gmail.com", "", "", "pashanwinsty1998@gmail.com")

# --- 5. Move the selection to match the new used range ---
$ws.Range("C4").Select() | Out-Null
